# Auto-generated Excel COM-interop script replicating the Kujata_Profits
# market-data refresh diff (scheduled runner price snapshot update).
# For each affected leve row, currentAveragePrice(NQ/HQ) and the derived
# LevePrice/LeveProfit columns are rewritten to the new snapshot values;
# a handful of cells that no longer have a value are cleared outright.

$wb = $excel.ActiveWorkbook

# ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1399.5
$ws.Range("J40").Value = 1300
$ws.Range("L40").Value = 1300
$ws.Range("N40").Value = -1650

# ALC row 51
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1000
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").ClearContents()

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 2849.923
$ws.Range("I107").Value = 3407.5715
$ws.Range("J107").Value = 2199.3333
$ws.Range("K107").Value = 3407.5715
$ws.Range("L107").Value = 2199.3333
$ws.Range("M107").Value = -1487.5715
$ws.Range("N107").Value = -6039.3333

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5469577
$ws.Range("I132").Value = 6175950
$ws.Range("K132").Value = 18527850
$ws.Range("M132").Value = -18525320

# ARM row 24
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 24900.5
$ws.Range("J24").Value = 24900.5
$ws.Range("L24").Value = 24900.5
$ws.Range("N24").Value = -25648.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5537.31
$ws.Range("I32").Value = 4492.4165
$ws.Range("J32").Value = 11023
$ws.Range("K32").Value = 4492.4165
$ws.Range("L32").Value = 11023
$ws.Range("M32").Value = -4205.4165
$ws.Range("N32").Value = -11597

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1093.625
$ws.Range("I74").Value = 469.375
$ws.Range("J74").Value = 2342.125
$ws.Range("K74").Value = 469.375
$ws.Range("L74").Value = 2342.125
$ws.Range("M74").Value = 404.625
$ws.Range("N74").Value = -4090.125

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1093.625
$ws.Range("I77").Value = 469.375
$ws.Range("J77").Value = 2342.125
$ws.Range("K77").Value = 2346.875
$ws.Range("L77").Value = 11710.625
$ws.Range("M77").Value = 2021.125
$ws.Range("N77").Value = -20446.625

# ARM row 100
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H100").Value = 24900.5
$ws.Range("J100").Value = 24900.5
$ws.Range("L100").Value = 24900.5
$ws.Range("N100").Value = -27064.5

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1914.3334
$ws.Range("I122").Value = 1412.3429
$ws.Range("J122").Value = 4424.2856
$ws.Range("K122").Value = 4237.028700000001
$ws.Range("L122").Value = 13272.8568
$ws.Range("M122").Value = -1787.028700000001
$ws.Range("N122").Value = -18172.8568

# ARM row 124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 19691.4
$ws.Range("J124").Value = 19691.4
$ws.Range("L124").Value = 19691.4
$ws.Range("N124").Value = -29511.4

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2387.425
$ws.Range("I132").Value = 2034.4482
$ws.Range("K132").Value = 6103.3446
$ws.Range("M132").Value = -3573.3446

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 58827776
$ws.Range("J86").Value = 3500
$ws.Range("L86").Value = 3500
$ws.Range("N86").Value = -5746

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 58827776
$ws.Range("J89").Value = 3500
$ws.Range("L89").Value = 17500
$ws.Range("N89").Value = -28732

# BSM row 97
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 45555.5
$ws.Range("I97").Value = 9999.5
$ws.Range("K97").Value = 9999.5
$ws.Range("M97").Value = -9008.5

# BSM row 130
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 41189.332
$ws.Range("J130").Value = 41189.332
$ws.Range("L130").Value = 41189.332
$ws.Range("N130").Value = -51229.332

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5730.75
$ws.Range("I134").Value = 914.2353000000001
$ws.Range("K134").Value = 2742.7059
$ws.Range("M134").Value = -207.7058999999999

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1681.8727
$ws.Range("I31").Value = 1681.8727
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1681.8727
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1386.8727
$ws.Range("N31").ClearContents()

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1681.8727
$ws.Range("I34").Value = 1681.8727
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1681.8727
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1479.8727
$ws.Range("N34").ClearContents()

# CRP row 52
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 32431.5
$ws.Range("J52").Value = 34776
$ws.Range("L52").Value = 34776
$ws.Range("N52").Value = -35364

# CRP row 86
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 3538303.8
$ws.Range("I86").Value = 5146889
$ws.Range("J86").Value = 53035.332
$ws.Range("K86").Value = 5146889
$ws.Range("L86").Value = 53035.332
$ws.Range("M86").Value = -5145766
$ws.Range("N86").Value = -55281.332

# CRP row 88
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

# CRP row 89
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H89").Value = 3538303.8
$ws.Range("I89").Value = 5146889
$ws.Range("J89").Value = 53035.332
$ws.Range("K89").Value = 25734445
$ws.Range("L89").Value = 265176.66
$ws.Range("M89").Value = -25728829
$ws.Range("N89").Value = -276408.66

# CRP row 91
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

# CUL row 3
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 12601.12
$ws.Range("I3").Value = 8066.077
$ws.Range("J3").Value = 17514.084
$ws.Range("K3").Value = 24198.231
$ws.Range("L3").Value = 52542.25199999999
$ws.Range("M3").Value = -24086.231
$ws.Range("N3").Value = -52766.25199999999

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 31297876
$ws.Range("J131").Value = 54642.82
$ws.Range("L131").Value = 163928.46
$ws.Range("N131").Value = -174008.46

# GSM row 18
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

# GSM row 68
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 20295
$ws.Range("J68").Value = 20295
$ws.Range("L68").Value = 20295
$ws.Range("N68").Value = -21917

# GSM row 71
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H71").Value = 20295
$ws.Range("J71").Value = 20295
$ws.Range("L71").Value = 60885
$ws.Range("N71").Value = -68997

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5134.1953
$ws.Range("I132").Value = 5470.0625
$ws.Range("J132").Value = 3940
$ws.Range("K132").Value = 16410.1875
$ws.Range("L132").Value = 11820
$ws.Range("M132").Value = -13880.1875
$ws.Range("N132").Value = -16880

# GSM row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 25346.143
$ws.Range("J136").Value = 25346.143
$ws.Range("L136").Value = 76038.429
$ws.Range("N136").Value = -81138.429

# LTW row 18
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1318.7
$ws.Range("I61").Value = 1318.7
$ws.Range("K61").Value = 1318.7
$ws.Range("M61").Value = -1116.7

# LTW row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1104.8889
$ws.Range("I93").Value = 1055.5
$ws.Range("J93").Value = 1500
$ws.Range("K93").Value = 1055.5
$ws.Range("L93").Value = 1500
$ws.Range("M93").Value = 192.5
$ws.Range("N93").Value = -3996

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1318.7
$ws.Range("I113").Value = 1318.7
$ws.Range("K113").Value = 1318.7
$ws.Range("M113").Value = 851.3

# LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 12321549
$ws.Range("I122").Value = 13494697
$ws.Range("K122").Value = 40484091
$ws.Range("M122").Value = -40481641

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1610.7368
$ws.Range("I136").Value = 1525.9333
$ws.Range("J136").Value = 1928.75
$ws.Range("K136").Value = 4577.7999
$ws.Range("L136").Value = 5786.25
$ws.Range("M136").Value = -2027.7999
$ws.Range("N136").Value = -10886.25

# WVR row 48
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H48").Value = 10032.5

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 46297216
$ws.Range("I126").Value = 50505970
$ws.Range("K126").Value = 151517910
$ws.Range("M126").Value = -151515440

# WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1712.659
$ws.Range("I136").Value = 674.85
$ws.Range("K136").Value = 2024.55
$ws.Range("M136").Value = 525.4499999999998
